# Apply weekly update: insert a new data row at row 694 (pushing existing
# rows 694-739 down to 695-740) and populate it with the new week's record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 694; this shifts rows 694-739 down to 695-740
$ws.Rows.Item(694).Insert()

# Populate the newly inserted row 694 with the new data point
$ws.Range("A694").Value = 9
$ws.Range("B694").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C694").Value = "Metropolitana"
$ws.Range("D694").Value = 44826
$ws.Range("D694").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E694").Value = 13
$ws.Range("F694").Value = 100112040
$ws.Range("G694").Value = "Cilantro"
$ws.Range("H694").Value = "Sin especificar"
$ws.Range("I694").Value = "Primera"
$ws.Range("J694").Value = 160
$ws.Range("K694").Value = 8000
$ws.Range("L694").Value = 10000
$ws.Range("M694").Value = 9000
$ws.Range("N694").Value = "`$/docena de atados"
$ws.Range("O694").Value = "Región Metropolitana"
$ws.Range("P694").Value = 3000
$ws.Range("Q694").Value = 3
$ws.Range("R694").Value = "Hortaliza"
